$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the data block (row 641),
# pushing the existing rows 641:677 down to 642:678.
$ws.Rows("641:641").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(641, 1).Value2 = 6
$ws.Cells.Item(641, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(641, 3).Value2 = "Metropolitana"
$ws.Cells.Item(641, 4).Value2 = 45267
$ws.Cells.Item(641, 5).Value2 = 13
$ws.Cells.Item(641, 6).Value2 = 100112032
$ws.Cells.Item(641, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(641, 8).Value2 = "Sin especificar"
$ws.Cells.Item(641, 9).Value2 = "Primera"
$ws.Cells.Item(641, 10).Value2 = 500
$ws.Cells.Item(641, 11).Value2 = 9000
$ws.Cells.Item(641, 12).Value2 = 10000
$ws.Cells.Item(641, 13).Value2 = 9540
$ws.Cells.Item(641, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(641, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(641, 16).Value2 = 191
$ws.Cells.Item(641, 17).Value2 = 50
$ws.Cells.Item(641, 18).Value2 = "Hortaliza"
